$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-10 (row 23)
$ws.Range("B23").Value = 6315
$ws.Range("C23").Value = 1003
$ws.Range("D23").Value = 5912487
$ws.Range("E23").Value = 936.2608076009501
$ws.Range("F23").Value = 8.356211393273849
$ws.Range("G23").Value = 4.370447450572312
$ws.Range("H23").Value = 26.68946260622387
